$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.700.75"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.294.32"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'302.29"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'96.03"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'34.79"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "'18.61"
$ws.Range("E12").Value = "  +5.24%  "
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("D14").Value = "'6.83"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "2.651.66"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "2.293.71"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "42.617.63"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "'12.75"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'6.00"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "'67.04"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "'236.00"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").Value = "'24.61"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").Value = "'167.22"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D31").Value = "'32.72"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'17.75"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("E35").Value = "  -6.85%  "
$ws.Range("D36").Value = "'2.36"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").Value = "'0.0684"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "'2.69"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("D42").Value = "1.993.75"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.18"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'10.19"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "'18.25"
$ws.Range("E46").Value = "  +6.00%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.92"
$ws.Range("E48").Value = "  +8.08%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'53.36"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "2.518.98"
$ws.Range("D51").Value = "'70.73"
$ws.Range("E51").Value = "  -1.92%  "
